$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 27: num_customers (C) and retention_rate (E)
$ws.Range("C27").Value = 52
$ws.Range("E27").Value = 0.02309058614564831

# Update row 37: num_customers (C) and cohort_size (D)
$ws.Range("C37").Value = 928
$ws.Range("D37").Value = 928
